$wb = $excel.ActiveWorkbook

# ALC row 105
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 48012.6
$ws.Range("J105").Value = 48012.6
$ws.Range("L105").Value = 48012.6
$ws.Range("N105").Value = -55000.6

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3859.4285
$ws.Range("I116").Value = 2581.5715
$ws.Range("J116").Value = 6415.143
$ws.Range("K116").Value = 2581.5715
$ws.Range("L116").Value = 6415.143
$ws.Range("M116").Value = 860.4285
$ws.Range("N116").Value = -13299.143

# ALC row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1167
$ws.Range("I127").Value = 613.5714
$ws.Range("J127").Value = 1268.9474
$ws.Range("K127").Value = 1840.7142
$ws.Range("L127").Value = 3806.8422
$ws.Range("M127").Value = 3119.2858
$ws.Range("N127").Value = -13726.8422

# ALC row 128
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H128").Value = 46608.668
$ws.Range("J128").Value = 46608.668
$ws.Range("L128").Value = 46608.668
$ws.Range("N128").Value = -56568.668

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10822.685
$ws.Range("I32").Value = 9821.156999999999
$ws.Range("K32").Value = 9821.156999999999
$ws.Range("M32").Value = -9534.156999999999

# ARM row 121
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 32658.6
$ws.Range("J121").Value = 32658.6
$ws.Range("L121").Value = 32658.6
$ws.Range("N121").Value = -36152.6

# ARM row 123
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 41775.5
$ws.Range("J123").Value = 41775.5
$ws.Range("L123").Value = 41775.5
$ws.Range("N123").Value = -51575.5

# ARM row 130
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H130").Value = 48421
$ws.Range("J130").Value = 48421
$ws.Range("L130").Value = 48421
$ws.Range("N130").Value = -58461

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 13890646
$ws.Range("I132").Value = 17242424
$ws.Range("J132").Value = 4708.5713
$ws.Range("K132").Value = 51727272
$ws.Range("L132").Value = 14125.7139
$ws.Range("M132").Value = -51724742
$ws.Range("N132").Value = -19185.7139

# ARM row 137
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# ARM row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 48539.8
$ws.Range("J139").Value = 48539.8
$ws.Range("L139").Value = 48539.8
$ws.Range("N139").Value = -58819.8

# BSM row 11
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 500
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3378.2273
$ws.Range("I105").Value = 2523.75
$ws.Range("J105").Value = 4403.6
$ws.Range("K105").Value = 2523.75
$ws.Range("L105").Value = 4403.6
$ws.Range("M105").Value = -776.75
$ws.Range("N105").Value = -7897.6

# BSM row 129
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999

# BSM row 130
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 49181
$ws.Range("J130").Value = 49181
$ws.Range("L130").Value = 49181
$ws.Range("N130").Value = -59221

# BSM row 132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 25502.695
$ws.Range("J132").Value = 25502.695
$ws.Range("L132").Value = 25502.695
$ws.Range("N132").Value = -35622.695

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2892.7903
$ws.Range("I134").Value = 1119.8182
$ws.Range("J134").Value = 3867.925
$ws.Range("K134").Value = 3359.4546
$ws.Range("L134").Value = 11603.775
$ws.Range("M134").Value = -824.4546
$ws.Range("N134").Value = -16673.775

# BSM row 140
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 43366
$ws.Range("J140").Value = 43366
$ws.Range("L140").Value = 43366
$ws.Range("N140").Value = -53726

# CRP row 20
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49772
$ws.Range("J20").Value = 49772
$ws.Range("L20").Value = 49772
$ws.Range("N20").Value = -50244

# CRP row 30
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H30").Value = 49772
$ws.Range("J30").Value = 49772
$ws.Range("L30").Value = 49772
$ws.Range("N30").Value = -49954

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3162.3
$ws.Range("I31").Value = 1120.8636
$ws.Range("J31").Value = 3738.0898
$ws.Range("K31").Value = 1120.8636
$ws.Range("L31").Value = 3738.0898
$ws.Range("M31").Value = -825.8635999999999
$ws.Range("N31").Value = -4328.0898

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3162.3
$ws.Range("I34").Value = 1120.8636
$ws.Range("J34").Value = 3738.0898
$ws.Range("K34").Value = 1120.8636
$ws.Range("L34").Value = 3738.0898
$ws.Range("M34").Value = -918.8635999999999
$ws.Range("N34").Value = -4142.0898

# CRP row 103
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 7697.364
$ws.Range("I103").Value = 4974.5557
$ws.Range("K103").Value = 4974.5557
$ws.Range("M103").Value = -3802.5557

# CRP row 116
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H116").Value = 49819.668
$ws.Range("J116").Value = 49819.668
$ws.Range("L116").Value = 49819.668
$ws.Range("N116").Value = -58997.668

# CRP row 128
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H128").Value = 49772
$ws.Range("J128").Value = 49772
$ws.Range("L128").Value = 49772
$ws.Range("N128").Value = -59732

# CRP row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 5055.857
$ws.Range("J141").Value = 5055.857
$ws.Range("L141").Value = 5055.857
$ws.Range("N141").Value = -15415.857

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4878.381
$ws.Range("I97").Value = 3369
$ws.Range("J97").Value = 8651.833000000001
$ws.Range("K97").Value = 3369
$ws.Range("L97").Value = 8651.833000000001
$ws.Range("M97").Value = -2873
$ws.Range("N97").Value = -9643.833000000001

# GSM row 104
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 43097.832
$ws.Range("J104").Value = 43097.832
$ws.Range("L104").Value = 43097.832
$ws.Range("N104").Value = -50085.832

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1255.2142
$ws.Range("I113").Value = 1126
$ws.Range("J113").Value = 1578.25
$ws.Range("K113").Value = 1126
$ws.Range("L113").Value = 1578.25
$ws.Range("M113").Value = 1044
$ws.Range("N113").Value = -5918.25

# GSM row 124
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 40845
$ws.Range("J124").Value = 40845
$ws.Range("L124").Value = 40845
$ws.Range("N124").Value = -50665

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 8010.25
$ws.Range("I126").Value = 9986.071
$ws.Range("J126").Value = 3400
$ws.Range("K126").Value = 29958.213
$ws.Range("L126").Value = 10200
$ws.Range("M126").Value = -27488.213
$ws.Range("N126").Value = -15140

# GSM row 130
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 53986
$ws.Range("J130").Value = 53986
$ws.Range("L130").Value = 53986
$ws.Range("N130").Value = -64026

# GSM row 131
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H131").Value = 32577.5
$ws.Range("J131").Value = 32577.5
$ws.Range("L131").Value = 32577.5
$ws.Range("N131").Value = -42657.5

# LTW row 69
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

# LTW row 72
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

# LTW row 98
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 47988
$ws.Range("J98").Value = 47988
$ws.Range("L98").Value = 47988
$ws.Range("N98").Value = -53978

# LTW row 127
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 50711
$ws.Range("J127").Value = 50711
$ws.Range("L127").Value = 50711
$ws.Range("N127").Value = -60631

# LTW row 128
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 41925
$ws.Range("J128").Value = 41925
$ws.Range("L128").Value = 41925
$ws.Range("N128").Value = -51885

# LTW row 129
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

# WVR row 32
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 23663.334
$ws.Range("J32").Value = 27995
$ws.Range("L32").Value = 27995
$ws.Range("N32").Value = -28629

# WVR row 103
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 49594
$ws.Range("J103").Value = 49594
$ws.Range("L103").Value = 49594
$ws.Range("N103").Value = -51938

# WVR row 128
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 42244
$ws.Range("J128").Value = 42244
$ws.Range("L128").Value = 42244
$ws.Range("N128").Value = -52204

# WVR row 131
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 48995.832
$ws.Range("J131").Value = 48995.832
$ws.Range("L131").Value = 48995.832
$ws.Range("N131").Value = -59075.832
